$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Ayati Arvind"

$wsContact = $wb.Worksheets.Item("Contact")
$wsContact.Range("E2").Value = "testexternalcontact@email.com"
